$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data rows 262:267 (Google mobility / MYS update) ----
$data = @(
    @(262, 44139, 1032, 3,  10339, 82, 27),
    @(263, 44140, 1009, 9,  10503, 78, 28),
    @(264, 44141, 1755, 3,  11530, 83, 32),
    @(265, 44142, 1168, 0,  11666, 87, 32),
    @(266, 44143, 852,  13, 11689, 94, 32),
    @(267, 44144, 972,  5,  11308, 86, 31)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Formula = "=D$r-F$r"
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Formula = "=F$r/D$r"
}
